$d = $word.ActiveDocument

# Remove the two duplicate "Associate Proffessor Sven Johansson" paragraphs
# (one plain, one bold) that precede the "Ass. Prof.Sven Johansson" details
# paragraph in the References section.

$start = $d.Paragraphs.Item(104).Range.Start
$end = $d.Paragraphs.Item(105).Range.End

$r = $d.Range($start, $end)
$r.Delete()
